$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("proveedores")

# "Proveedor Alfa" (row 2) estado changes from "Inactivo" to "Activo"
$ws.Range("H2").Value = "Activo"

# New proveedor added as row 8: "Gisela porfiri proveedor"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "Gisela porfiri proveedor"

# cuit/telefono look numeric ("123123" / "12345") but must be stored as text,
# matching how the rest of the sheet stores such values.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "123123"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "gisela2@email.com"

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "12345"
$ws.Range("F8").ClearFormats()

$ws.Range("G8").Value = "del valle 462"
$ws.Range("H8").Value = "Inactivo"
